$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.613.62"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").Value = "3.129.22"
$ws.Range("E3").Value = "  -2.33%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.89"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.58"
$ws.Range("E6").Value = "  -5.45%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.124.06"
$ws.Range("E8").Value = "  -2.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.31"
$ws.Range("E11").Value = "  -4.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("E13").Value = "  -3.01%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.32"
$ws.Range("E14").Value = "  -5.48%  "

$ws.Range("D15").Value = "3.635.86"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("E16").Value = "  +2.30%  "

$ws.Range("D17").Value = "63.643.56"
$ws.Range("E17").Value = "  -2.24%  "

$ws.Range("D18").Value = "3.123.21"
$ws.Range("E18").Value = "  -2.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  -3.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.05"
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  -2.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.709"
$ws.Range("E22").Value = "  -2.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.89"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.59"
$ws.Range("E24").Value = "  -3.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.99"
$ws.Range("E25").Value = "  -2.75%  "

$ws.Range("E26").Value = "  +0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  -5.49%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("E28").Value = "  -3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.120"
$ws.Range("E29").Value = "  -7.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.06"
$ws.Range("E31").Value = "  -10.13%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.68"
$ws.Range("E33").Value = "  -2.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.12"
$ws.Range("E34").Value = "  -3.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  +1.00%  "

$ws.Range("D36").Value = "0.0₃0775"
$ws.Range("E36").Value = "  +5.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.97"
$ws.Range("E37").Value = "  -3.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.55"
$ws.Range("E38").Value = "  -4.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "456.91"
$ws.Range("E39").Value = "  -4.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("E40").Value = "  -9.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("E41").Value = "  -3.32%  "

$ws.Range("E42").Value = "  -6.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.29"
$ws.Range("E43").Value = "  -3.28%  "

$ws.Range("D44").Value = "2.852.58"
$ws.Range("E44").Value = "  -2.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.267"
$ws.Range("E45").Value = "  -4.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("E46").Value = "  -7.15%  "

$ws.Range("E47").Value = "  +1.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.24"
$ws.Range("E49").Value = "  -5.19%  "

$ws.Range("E50").Value = "  -2.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.03"
$ws.Range("E51").Value = "  -1.63%  "
